$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sami/bronze -> samis/platinum
$ws.Cells.Item(2, 1).Value = "samis"
$ws.Cells.Item(2, 3).Value = "platinum"

# Row 3: samis/platinum -> ds/bronze
$ws.Cells.Item(3, 1).Value = "ds"
$ws.Cells.Item(3, 3).Value = "bronze"

# Row 4: sami/gold -> AA/bronze
$ws.Cells.Item(4, 1).Value = "AA"
$ws.Cells.Item(4, 3).Value = "bronze"

# Row 5: sami/bronze/sgs -- unchanged

# Row 6 (new): sami / TRUE / Master / sgs
$ws.Cells.Item(6, 1).Value = "sami"
$ws.Cells.Item(6, 2).Value = $true
$ws.Cells.Item(6, 3).Value = "Master"
$ws.Cells.Item(6, 4).Value = "sgs"
